# PVDI_11_Inventor_Canonical_Test_Cases.xlsx
# "Extended canonical name assignment to pregrant data"
#
# Content change: the 4th column header ("patent_date") is renamed to
# "doc_date" on all three worksheets, reflecting that the date column now
# also covers pregrant-publication documents, not just granted patents.
#
# Cosmetic/view-state changes that accompanied the save: the active sheet
# moved from "Tie_Resolution_Case" (3rd tab) to "Multiple_Names_Case" (1st
# tab), and each sheet's last-used-cell selection was updated.

$wb = $excel.ActiveWorkbook

$wsMultiple = $wb.Worksheets.Item("Multiple_Names_Case")
$wsSingle   = $wb.Worksheets.Item("Single_Data_Point_Case")
$wsTie      = $wb.Worksheets.Item("Tie_Resolution_Case")

# Rename the date column header on every sheet.
$wsMultiple.Range("D1").Value = "doc_date"
$wsSingle.Range("D1").Value = "doc_date"
$wsTie.Range("D1").Value = "doc_date"

# Update the remembered selection on the sheets that stay inactive first …
$wsSingle.Range("D7").Select() | Out-Null
$wsTie.Range("D2").Select() | Out-Null

# … then activate Multiple_Names_Case and leave it as the selected tab with
# its own selection, matching the final saved state.
$wsMultiple.Activate() | Out-Null
$wsMultiple.Range("D2").Select() | Out-Null
